$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $CellRef, $TextValue)
    $range = $Sheet.Range($CellRef)
    $range.Formula = "'" + $TextValue
    $range.Style = "Normal"
}

Set-TextValue $ws 'D2' '51.673.81'
Set-TextValue $ws 'E2' '  +0.86%  '
Set-TextValue $ws 'D3' '2.990.50'
Set-TextValue $ws 'E3' '  +2.38%  '
Set-TextValue $ws 'E4' '  +0.20%  '
Set-TextValue $ws 'D5' '384.29'
Set-TextValue $ws 'E5' '  +2.90%  '
Set-TextValue $ws 'D6' '104.60'
Set-TextValue $ws 'E6' '  +2.12%  '
Set-TextValue $ws 'E7' '  +0.68%  '
Set-TextValue $ws 'E8' '  +0.08%  '
Set-TextValue $ws 'E9' '  +1.37%  '
Set-TextValue $ws 'D10' '37.26'
Set-TextValue $ws 'E10' '  +0.64%  '
Set-TextValue $ws 'E11' '  +0.24%  '
Set-TextValue $ws 'D12' '0.0851'
Set-TextValue $ws 'E12' '  +1.82%  '
Set-TextValue $ws 'D13' '3.458.54'
Set-TextValue $ws 'E13' '  +2.77%  '
Set-TextValue $ws 'D14' '18.41'
Set-TextValue $ws 'E14' '  +0.37%  '
Set-TextValue $ws 'D15' '7.62'
Set-TextValue $ws 'E15' '  +2.42%  '
Set-TextValue $ws 'D16' '2.987.54'
Set-TextValue $ws 'E16' '  +2.66%  '
Set-TextValue $ws 'E17' '  +8.82%  '
Set-TextValue $ws 'D18' '51.591.56'
Set-TextValue $ws 'E18' '  +0.91%  '
Set-TextValue $ws 'E19' '  +0.38%  '
Set-TextValue $ws 'E20' '  +2.88%  '
Set-TextValue $ws 'D21' '12.91'
Set-TextValue $ws 'E21' '  -0.04%  '
Set-TextValue $ws 'D22' '0.0₃0967'
Set-TextValue $ws 'E22' '  +2.27%  '
Set-TextValue $ws 'D23' '69.28'
Set-TextValue $ws 'E23' '  +1.41%  '
Set-TextValue $ws 'D24' '263.42'
Set-TextValue $ws 'E24' '  +1.29%  '
Set-TextValue $ws 'E25' '  +8.38%  '
Set-TextValue $ws 'D26' '8.45'
Set-TextValue $ws 'E26' '  +18.17%  '
Set-TextValue $ws 'D27' '7.80'
Set-TextValue $ws 'E27' '  +17.85%  '
Set-TextValue $ws 'D28' '0.116'
Set-TextValue $ws 'E28' '  +14.15%  '
Set-TextValue $ws 'E29' '  +0.14%  '
Set-TextValue $ws 'D30' '26.10'
Set-TextValue $ws 'E30' '  +1.39%  '
Set-TextValue $ws 'E31' '  -0.13%  '
Set-TextValue $ws 'E32' '  +0.30%  '
Set-TextValue $ws 'E33' '  +1.17%  '
Set-TextValue $ws 'D34' '51.09'
Set-TextValue $ws 'E34' '  -0.46%  '
Set-TextValue $ws 'E35' '  -1.97%  '
Set-TextValue $ws 'E36' '  +6.92%  '
Set-TextValue $ws 'E37' '  +0.06%  '
Set-TextValue $ws 'D38' '3.05'
Set-TextValue $ws 'E38' '  +1.80%  '
Set-TextValue $ws 'E39' '  -0.17%  '
Set-TextValue $ws 'D40' '2.60'
Set-TextValue $ws 'E40' '  +0.71%  '
Set-TextValue $ws 'E41' '  +3.17%  '
Set-TextValue $ws 'E42' '  -0.25%  '
Set-TextValue $ws 'D43' '122.51'
Set-TextValue $ws 'E43' '  +2.33%  '
Set-TextValue $ws 'D44' '21.85'
Set-TextValue $ws 'E44' '  -1.15%  '
Set-TextValue $ws 'E45' '  +15.80%  '
Set-TextValue $ws 'E46' '  -1.83%  '
Set-TextValue $ws 'E47' '  +2.75%  '
Set-TextValue $ws 'D48' '3.32'
Set-TextValue $ws 'E48' '  +4.79%  '
Set-TextValue $ws 'D49' '2.039.74'
Set-TextValue $ws 'E49' '  +1.00%  '
Set-TextValue $ws 'D50' '0.0334'
Set-TextValue $ws 'E50' '  +7.24%  '
Set-TextValue $ws 'E51' '  +2.06%  '
